$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix correction sensitivity: 100Ω Resistor -> 110Ω Resistor
$ws.Range("A9").Value = "110Ω Resistor"

# Update selection (cursor moved to K16 as part of editing)
$ws.Range("K16").Select()
